$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(44, 8).Value = 0  # H44: 15000 -> 0
$ws.Cells.Item(44, 10).Value = 0  # J44: 15000 -> 0
$ws.Cells.Item(44, 12).Value = 0  # L44: 15000 -> 0
$ws.Cells.Item(44, 14).ClearContents()  # N44: -15924 -> (blank)
$ws.Cells.Item(51, 8).Value = 4000  # H51: 3999.5 -> 4000
$ws.Cells.Item(51, 9).Value = 4000  # I51: 3999.5 -> 4000
$ws.Cells.Item(51, 11).Value = 4000  # K51: 3999.5 -> 4000
$ws.Cells.Item(51, 13).Value = -3516  # M51: -3515.5 -> -3516
$ws.Cells.Item(52, 8).Value = 2000  # H52: 1799.8 -> 2000
$ws.Cells.Item(52, 10).Value = 2000  # J52: 1799.8 -> 2000
$ws.Cells.Item(52, 12).Value = 6000  # L52: 5399.4 -> 6000
$ws.Cells.Item(52, 14).Value = -6320  # N52: -5719.4 -> -6320
$ws.Cells.Item(55, 8).Value = 443.66666  # H55: 512.8333 -> 443.66666
$ws.Cells.Item(55, 9).Value = 40  # I55: 245 -> 40
$ws.Cells.Item(55, 10).Value = 524.4  # J55: 646.75 -> 524.4
$ws.Cells.Item(55, 11).Value = 40  # K55: 245 -> 40
$ws.Cells.Item(55, 12).Value = 524.4  # L55: 646.75 -> 524.4
$ws.Cells.Item(55, 13).Value = 174  # M55: -31 -> 174
$ws.Cells.Item(55, 14).Value = -952.4  # N55: -1074.75 -> -952.4
$ws.Cells.Item(62, 8).Value = 3897.1667  # H62: 3680.9443 -> 3897.1667
$ws.Cells.Item(62, 9).Value = 3576.6  # I62: 3453.5625 -> 3576.6
$ws.Cells.Item(62, 11).Value = 3576.6  # K62: 3453.5625 -> 3576.6
$ws.Cells.Item(62, 13).Value = -2952.6  # M62: -2829.5625 -> -2952.6
$ws.Cells.Item(65, 8).Value = 3897.1667  # H65: 3680.9443 -> 3897.1667
$ws.Cells.Item(65, 9).Value = 3576.6  # I65: 3453.5625 -> 3576.6
$ws.Cells.Item(65, 11).Value = 17883  # K65: 17267.8125 -> 17883
$ws.Cells.Item(65, 13).Value = -14763  # M65: -14147.8125 -> -14763
$ws.Cells.Item(76, 8).Value = 5665.75  # H76: 4936.3335 -> 5665.75
$ws.Cells.Item(76, 9).Value = 4919.6665  # I76: 4342.8 -> 4919.6665
$ws.Cells.Item(76, 11).Value = 4919.6665  # K76: 4342.8 -> 4919.6665
$ws.Cells.Item(76, 13).Value = -4604.6665  # M76: -4027.8 -> -4604.6665
$ws.Cells.Item(79, 8).Value = 5665.75  # H79: 4936.3335 -> 5665.75
$ws.Cells.Item(79, 9).Value = 4919.6665  # I79: 4342.8 -> 4919.6665
$ws.Cells.Item(79, 11).Value = 4919.6665  # K79: 4342.8 -> 4919.6665
$ws.Cells.Item(79, 13).Value = -3827.6665  # M79: -3250.8 -> -3827.6665
$ws.Cells.Item(86, 8).Value = 15000  # H86: 8200 -> 15000
$ws.Cells.Item(86, 9).Value = 15000  # I86: 8200 -> 15000
$ws.Cells.Item(86, 11).Value = 15000  # K86: 8200 -> 15000
$ws.Cells.Item(86, 13).Value = -13877  # M86: -7077 -> -13877
$ws.Cells.Item(89, 8).Value = 15000  # H89: 8200 -> 15000
$ws.Cells.Item(89, 9).Value = 15000  # I89: 8200 -> 15000
$ws.Cells.Item(89, 11).Value = 75000  # K89: 41000 -> 75000
$ws.Cells.Item(89, 13).Value = -69384  # M89: -35384 -> -69384

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1692.8572  # H32: 1887.1831 -> 1692.8572
$ws.Cells.Item(32, 9).Value = 1572.4637  # I32: 1753.4783 -> 1572.4637
$ws.Cells.Item(32, 10).Value = 10000  # J32: 6500 -> 10000
$ws.Cells.Item(32, 11).Value = 1572.4637  # K32: 1753.4783 -> 1572.4637
$ws.Cells.Item(32, 12).Value = 10000  # L32: 6500 -> 10000
$ws.Cells.Item(32, 13).Value = -1285.4637  # M32: -1466.4783 -> -1285.4637
$ws.Cells.Item(32, 14).Value = -10574  # N32: -7074 -> -10574
$ws.Cells.Item(45, 8).Value = 1275.8334  # H45: 1518.125 -> 1275.8334
$ws.Cells.Item(45, 9).Value = 1283  # I45: 1400.8334 -> 1283
$ws.Cells.Item(45, 10).Value = 1240  # J45: 1870 -> 1240
$ws.Cells.Item(45, 11).Value = 1283  # K45: 1400.8334 -> 1283
$ws.Cells.Item(45, 12).Value = 1240  # L45: 1870 -> 1240
$ws.Cells.Item(45, 13).Value = -906  # M45: -1023.8334 -> -906
$ws.Cells.Item(45, 14).Value = -1994  # N45: -2624 -> -1994
$ws.Cells.Item(55, 8).Value = 28333.334  # H55: 27498.125 -> 28333.334
$ws.Cells.Item(55, 10).Value = 28333.334  # J55: 27498.125 -> 28333.334
$ws.Cells.Item(55, 12).Value = 28333.334  # L55: 27498.125 -> 28333.334
$ws.Cells.Item(55, 14).Value = -28963.334  # N55: -28128.125 -> -28963.334
$ws.Cells.Item(63, 8).Value = 4451  # H63: 4551.2 -> 4451
$ws.Cells.Item(63, 9).Value = 3999  # I63: 0 -> 3999
$ws.Cells.Item(63, 10).Value = 4564  # J63: 4551.2 -> 4564
$ws.Cells.Item(63, 11).Value = 3999  # K63: 0 -> 3999
$ws.Cells.Item(63, 12).Value = 4564  # L63: 4551.2 -> 4564
$ws.Cells.Item(63, 13).Value = -3313  # M63: (blank) -> -3313
$ws.Cells.Item(63, 14).Value = -5936  # N63: -5923.2 -> -5936
$ws.Cells.Item(66, 8).Value = 4451  # H66: 4551.2 -> 4451
$ws.Cells.Item(66, 9).Value = 3999  # I66: 0 -> 3999
$ws.Cells.Item(66, 10).Value = 4564  # J66: 4551.2 -> 4564
$ws.Cells.Item(66, 11).Value = 19995  # K66: 0 -> 19995
$ws.Cells.Item(66, 12).Value = 22820  # L66: 22756 -> 22820
$ws.Cells.Item(66, 13).Value = -16563  # M66: (blank) -> -16563
$ws.Cells.Item(66, 14).Value = -29684  # N66: -29620 -> -29684
$ws.Cells.Item(74, 8).Value = 890.3333  # H74: 910.6 -> 890.3333
$ws.Cells.Item(74, 9).Value = 890.3333  # I74: 910.6 -> 890.3333
$ws.Cells.Item(74, 11).Value = 890.3333  # K74: 910.6 -> 890.3333
$ws.Cells.Item(74, 13).Value = -16.33330000000001  # M74: -36.60000000000002 -> -16.33330000000001
$ws.Cells.Item(77, 8).Value = 890.3333  # H77: 910.6 -> 890.3333
$ws.Cells.Item(77, 9).Value = 890.3333  # I77: 910.6 -> 890.3333
$ws.Cells.Item(77, 11).Value = 4451.6665  # K77: 4553 -> 4451.6665
$ws.Cells.Item(77, 13).Value = -83.66650000000027  # M77: -185 -> -83.66650000000027
$ws.Cells.Item(98, 8).Value = 14901.667  # H98: 18927.5 -> 14901.667
$ws.Cells.Item(98, 10).Value = 14901.667  # J98: 18927.5 -> 14901.667
$ws.Cells.Item(98, 12).Value = 14901.667  # L98: 18927.5 -> 14901.667
$ws.Cells.Item(98, 14).Value = -20891.667  # N98: -24917.5 -> -20891.667
$ws.Cells.Item(105, 8).Value = 0  # H105: 20000 -> 0
$ws.Cells.Item(105, 10).Value = 0  # J105: 20000 -> 0
$ws.Cells.Item(105, 12).Value = 0  # L105: 20000 -> 0
$ws.Cells.Item(105, 14).ClearContents()  # N105: -26988 -> (blank)
$ws.Cells.Item(125, 8).Value = 22000  # H125: 0 -> 22000
$ws.Cells.Item(125, 10).Value = 22000  # J125: 0 -> 22000
$ws.Cells.Item(125, 12).Value = 22000  # L125: 0 -> 22000
$ws.Cells.Item(125, 14).Value = -31840  # N125: (blank) -> -31840
$ws.Cells.Item(132, 8).Value = 2290.5881  # H132: 2571.0715 -> 2290.5881
$ws.Cells.Item(132, 9).Value = 2290.5881  # I132: 2571.0715 -> 2290.5881
$ws.Cells.Item(132, 11).Value = 6871.7643  # K132: 7713.2145 -> 6871.7643
$ws.Cells.Item(132, 13).Value = -4341.7643  # M132: -5183.2145 -> -4341.7643

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3754.5  # H86: 2965.6428 -> 3754.5
$ws.Cells.Item(86, 9).Value = 2174.1667  # I86: 1701.9 -> 2174.1667
$ws.Cells.Item(86, 11).Value = 2174.1667  # K86: 1701.9 -> 2174.1667
$ws.Cells.Item(86, 13).Value = -1051.1667  # M86: -578.9000000000001 -> -1051.1667
$ws.Cells.Item(89, 8).Value = 3754.5  # H89: 2965.6428 -> 3754.5
$ws.Cells.Item(89, 9).Value = 2174.1667  # I89: 1701.9 -> 2174.1667
$ws.Cells.Item(89, 11).Value = 10870.8335  # K89: 8509.5 -> 10870.8335
$ws.Cells.Item(89, 13).Value = -5254.833500000001  # M89: -2893.5 -> -5254.833500000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1000  # H16: 0 -> 1000
$ws.Cells.Item(16, 9).Value = 1000  # I16: 0 -> 1000
$ws.Cells.Item(16, 11).Value = 1000  # K16: 0 -> 1000
$ws.Cells.Item(16, 13).Value = -713  # M16: (blank) -> -713
$ws.Cells.Item(31, 8).Value = 3310.6155  # H31: 3514.8 -> 3310.6155
$ws.Cells.Item(31, 9).Value = 2136.125  # I31: 2199.8333 -> 2136.125
$ws.Cells.Item(31, 10).Value = 5189.8  # J31: 5487.25 -> 5189.8
$ws.Cells.Item(31, 11).Value = 2136.125  # K31: 2199.8333 -> 2136.125
$ws.Cells.Item(31, 12).Value = 5189.8  # L31: 5487.25 -> 5189.8
$ws.Cells.Item(31, 13).Value = -1841.125  # M31: -1904.8333 -> -1841.125
$ws.Cells.Item(31, 14).Value = -5779.8  # N31: -6077.25 -> -5779.8
$ws.Cells.Item(34, 8).Value = 3310.6155  # H34: 3514.8 -> 3310.6155
$ws.Cells.Item(34, 9).Value = 2136.125  # I34: 2199.8333 -> 2136.125
$ws.Cells.Item(34, 10).Value = 5189.8  # J34: 5487.25 -> 5189.8
$ws.Cells.Item(34, 11).Value = 2136.125  # K34: 2199.8333 -> 2136.125
$ws.Cells.Item(34, 12).Value = 5189.8  # L34: 5487.25 -> 5189.8
$ws.Cells.Item(34, 13).Value = -1934.125  # M34: -1997.8333 -> -1934.125
$ws.Cells.Item(34, 14).Value = -5593.8  # N34: -5891.25 -> -5593.8
$ws.Cells.Item(70, 8).Value = 24990  # H70: 25000 -> 24990
$ws.Cells.Item(70, 10).Value = 24990  # J70: 25000 -> 24990
$ws.Cells.Item(70, 12).Value = 24990  # L70: 25000 -> 24990
$ws.Cells.Item(70, 14).Value = -25620  # N70: -25630 -> -25620
$ws.Cells.Item(73, 8).Value = 24990  # H73: 25000 -> 24990
$ws.Cells.Item(73, 10).Value = 24990  # J73: 25000 -> 24990
$ws.Cells.Item(73, 12).Value = 24990  # L73: 25000 -> 24990
$ws.Cells.Item(73, 14).Value = -27174  # N73: -27184 -> -27174
$ws.Cells.Item(86, 8).Value = 7878.769  # H86: 7879.154 -> 7878.769
$ws.Cells.Item(86, 9).Value = 7292  # I86: 7292.5557 -> 7292
$ws.Cells.Item(86, 11).Value = 7292  # K86: 7292.5557 -> 7292
$ws.Cells.Item(86, 13).Value = -6169  # M86: -6169.5557 -> -6169
$ws.Cells.Item(89, 8).Value = 7878.769  # H89: 7879.154 -> 7878.769
$ws.Cells.Item(89, 9).Value = 7292  # I89: 7292.5557 -> 7292
$ws.Cells.Item(89, 11).Value = 36460  # K89: 36462.7785 -> 36460
$ws.Cells.Item(89, 13).Value = -30844  # M89: -30846.7785 -> -30844
$ws.Cells.Item(107, 8).Value = 278.875  # H107: 281.69232 -> 278.875
$ws.Cells.Item(107, 9).Value = 277.26666  # I107: 281.7647 -> 277.26666
$ws.Cells.Item(107, 11).Value = 277.26666  # K107: 281.7647 -> 277.26666
$ws.Cells.Item(107, 13).Value = 1642.73334  # M107: 1638.2353 -> 1642.73334
$ws.Cells.Item(113, 8).Value = 1000  # H113: 0 -> 1000
$ws.Cells.Item(113, 9).Value = 1000  # I113: 0 -> 1000
$ws.Cells.Item(113, 11).Value = 1000  # K113: 0 -> 1000
$ws.Cells.Item(113, 13).Value = 1170  # M113: (blank) -> 1170

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 1001500  # H11: 667766.3 -> 1001500
$ws.Cells.Item(11, 10).Value = 0  # J11: 299 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 897 -> 0
$ws.Cells.Item(11, 14).ClearContents()  # N11: -1177 -> (blank)
$ws.Cells.Item(80, 8).Value = 8966  # H80: 10238 -> 8966
$ws.Cells.Item(80, 10).Value = 8499.25  # J80: 10463.667 -> 8499.25
$ws.Cells.Item(80, 12).Value = 25497.75  # L80: 31391.001 -> 25497.75
$ws.Cells.Item(80, 14).Value = -27369.75  # N80: -33263.001 -> -27369.75
$ws.Cells.Item(83, 8).Value = 8966  # H83: 10238 -> 8966
$ws.Cells.Item(83, 10).Value = 8499.25  # J83: 10463.667 -> 8499.25
$ws.Cells.Item(83, 12).Value = 76493.25  # L83: 94173.003 -> 76493.25
$ws.Cells.Item(83, 14).Value = -85853.25  # N83: -103533.003 -> -85853.25
$ws.Cells.Item(122, 8).Value = 1497.8  # H122: 1623.75 -> 1497.8
$ws.Cells.Item(122, 10).Value = 1622.25  # J122: 1831.6666 -> 1622.25
$ws.Cells.Item(122, 12).Value = 14600.25  # L122: 16484.9994 -> 14600.25
$ws.Cells.Item(122, 14).Value = -19500.25  # N122: -21384.9994 -> -19500.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 641.875  # H97: 675.93335 -> 641.875
$ws.Cells.Item(97, 9).Value = 612.4167  # I97: 625.8182 -> 612.4167
$ws.Cells.Item(97, 10).Value = 730.25  # J97: 813.75 -> 730.25
$ws.Cells.Item(97, 11).Value = 612.4167  # K97: 625.8182 -> 612.4167
$ws.Cells.Item(97, 12).Value = 730.25  # L97: 813.75 -> 730.25
$ws.Cells.Item(97, 13).Value = -116.4167  # M97: -129.8182 -> -116.4167
$ws.Cells.Item(97, 14).Value = -1722.25  # N97: -1805.75 -> -1722.25
$ws.Cells.Item(123, 8).Value = 21000  # H123: 0 -> 21000
$ws.Cells.Item(123, 10).Value = 21000  # J123: 0 -> 21000
$ws.Cells.Item(123, 12).Value = 21000  # L123: 0 -> 21000
$ws.Cells.Item(123, 14).Value = -25900  # N123: (blank) -> -25900
$ws.Cells.Item(126, 8).Value = 5298.3  # H126: 8354.272000000001 -> 5298.3
$ws.Cells.Item(126, 9).Value = 6347.875  # I126: 12499.714 -> 6347.875
$ws.Cells.Item(126, 10).Value = 1100  # J126: 1099.75 -> 1100
$ws.Cells.Item(126, 11).Value = 19043.625  # K126: 37499.142 -> 19043.625
$ws.Cells.Item(126, 12).Value = 3300  # L126: 3299.25 -> 3300
$ws.Cells.Item(126, 13).Value = -16573.625  # M126: -35029.142 -> -16573.625
$ws.Cells.Item(126, 14).Value = -8240  # N126: -8239.25 -> -8240

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 0  # H16: 1197.5 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 895 -> 0
$ws.Cells.Item(16, 10).Value = 0  # J16: 1500 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 895 -> 0
$ws.Cells.Item(16, 12).Value = 0  # L16: 1500 -> 0
$ws.Cells.Item(16, 13).ClearContents()  # M16: -725 -> (blank)
$ws.Cells.Item(16, 14).ClearContents()  # N16: -1840 -> (blank)
$ws.Cells.Item(46, 8).Value = 3290.5908  # H46: 3385.9546 -> 3290.5908
$ws.Cells.Item(46, 9).Value = 2489.2  # I46: 2554.5557 -> 2489.2
$ws.Cells.Item(46, 10).Value = 3958.4167  # J46: 3961.5386 -> 3958.4167
$ws.Cells.Item(46, 11).Value = 2489.2  # K46: 2554.5557 -> 2489.2
$ws.Cells.Item(46, 12).Value = 3958.4167  # L46: 3961.5386 -> 3958.4167
$ws.Cells.Item(46, 13).Value = -2301.2  # M46: -2366.5557 -> -2301.2
$ws.Cells.Item(46, 14).Value = -4334.4167  # N46: -4337.5386 -> -4334.4167
$ws.Cells.Item(100, 8).Value = 3399.8  # H100: 3666.5 -> 3399.8
$ws.Cells.Item(100, 10).Value = 4999.6665  # J100: 4999.75 -> 4999.6665
$ws.Cells.Item(100, 12).Value = 4999.6665  # L100: 4999.75 -> 4999.6665
$ws.Cells.Item(100, 14).Value = -6081.6665  # N100: -6081.75 -> -6081.6665
$ws.Cells.Item(127, 8).Value = 76998.2  # H127: 77498 -> 76998.2
$ws.Cells.Item(127, 10).Value = 76998.2  # J127: 77498 -> 76998.2
$ws.Cells.Item(127, 12).Value = 76998.2  # L127: 77498 -> 76998.2
$ws.Cells.Item(127, 14).Value = -86918.2  # N127: -87418 -> -86918.2
$ws.Cells.Item(136, 8).Value = 3009.1667  # H136: 2919.0908 -> 3009.1667
$ws.Cells.Item(136, 10).Value = 4002.5  # J136: 4005 -> 4002.5
$ws.Cells.Item(136, 12).Value = 12007.5  # L136: 12015 -> 12007.5
$ws.Cells.Item(136, 14).Value = -17107.5  # N136: -17115 -> -17107.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 35161.5  # H41: 32388.8 -> 35161.5
$ws.Cells.Item(41, 9).Value = 50000  # I41: 0 -> 50000
$ws.Cells.Item(41, 10).Value = 32193.8  # J41: 32388.8 -> 32193.8
$ws.Cells.Item(41, 11).Value = 50000  # K41: 0 -> 50000
$ws.Cells.Item(41, 12).Value = 32193.8  # L41: 32388.8 -> 32193.8
$ws.Cells.Item(41, 13).Value = -49610  # M41: (blank) -> -49610
$ws.Cells.Item(41, 14).Value = -32973.8  # N41: -33168.8 -> -32973.8
